$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 42, pushing existing rows 42-118 down to 43-119
# (this reproduces the "new weekly observation" commit: the whole
# history table shifts down by one row and a brand-new row 42 is
# populated with the latest price observation).
$ws.Rows(42).Insert()

# Populate the new row 42 with the latest observation.
$ws.Cells.Item(42, 1).Value = 9
$ws.Cells.Item(42, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(42, 3).Value = "Metropolitana"
$ws.Cells.Item(42, 4).Value = 45203
$ws.Cells.Item(42, 5).Value = 13
$ws.Cells.Item(42, 6).Value = 100112029
$ws.Cells.Item(42, 7).Value = "Orégano"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 16
$ws.Cells.Item(42, 11).Value = 21000
$ws.Cells.Item(42, 12).Value = 21000
$ws.Cells.Item(42, 13).Value = 21000
$ws.Cells.Item(42, 14).Value = "$/docena de atados"
$ws.Cells.Item(42, 15).Value = "Región Metropolitana"
$ws.Cells.Item(42, 16).Value = 7000
$ws.Cells.Item(42, 17).Value = 3
$ws.Cells.Item(42, 18).Value = "Hortaliza"
